$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 26, pushing the existing rows 26-80 down to 27-81.
$ws.Rows.Item(26).Insert()

# Populate the newly inserted row 26 with the new weekly record.
$ws.Cells.Item(26, 1).Value = 5
$ws.Cells.Item(26, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(26, 3).Value = "Maule"
$ws.Cells.Item(26, 4).Value = 44883
$ws.Cells.Item(26, 5).Value = 7
$ws.Cells.Item(26, 6).Value = "Fruta"
$ws.Cells.Item(26, 7).Value = 100101
$ws.Cells.Item(26, 8).Value = "Berries"
$ws.Cells.Item(26, 9).Value = 100101001
$ws.Cells.Item(26, 10).Value = "Arándano (blue)"
$ws.Cells.Item(26, 11).Value = "Sin especificar"
$ws.Cells.Item(26, 12).Value = "Primera"
$ws.Cells.Item(26, 13).Value = 250
$ws.Cells.Item(26, 14).Value = 5600
$ws.Cells.Item(26, 15).Value = 6000
$ws.Cells.Item(26, 16).Value = 5760
$ws.Cells.Item(26, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(26, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(26, 19).Value = 2880
$ws.Cells.Item(26, 20).Value = 2
